# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Each entry is (row, column, new text value). Column D sometimes holds a
# value that is lexically a plain number ("212.56") rather than a
# thousands-grouped price ("27.906.80"); those must be forced to Text so
# Excel's COM layer doesn't silently coerce the cell to a Number (which
# would also reformat "151.00" -> 151, "6.90" -> 6.9, etc. and lose the
# original trailing zeros). We flip NumberFormat to "@" for the write and
# then restore the "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = 4; Value = "27.906.80"; ForceText = $false }
    @{ Row = 2; Col = 5; Value = "  +1.27%  "; ForceText = $false }
    @{ Row = 3; Col = 4; Value = "1.640.66"; ForceText = $false }
    @{ Row = 3; Col = 5; Value = "  +1.14%  "; ForceText = $false }
    @{ Row = 5; Col = 4; Value = "212.56"; ForceText = $true }
    @{ Row = 5; Col = 5; Value = "  +0.39%  "; ForceText = $false }
    @{ Row = 6; Col = 5; Value = "  +0.97%  "; ForceText = $false }
    @{ Row = 7; Col = 5; Value = "  -0.27%  "; ForceText = $false }
    @{ Row = 8; Col = 4; Value = "23.47"; ForceText = $true }
    @{ Row = 8; Col = 5; Value = "  +1.81%  "; ForceText = $false }
    @{ Row = 9; Col = 5; Value = "  +0.94%  "; ForceText = $false }
    @{ Row = 10; Col = 5; Value = "  +0.18%  "; ForceText = $false }
    @{ Row = 11; Col = 5; Value = "  -1.98%  "; ForceText = $false }
    @{ Row = 12; Col = 4; Value = "1.871.53"; ForceText = $false }
    @{ Row = 12; Col = 5; Value = "  +1.02%  "; ForceText = $false }
    @{ Row = 13; Col = 4; Value = "1.639.28"; ForceText = $false }
    @{ Row = 13; Col = 5; Value = "  +0.97%  "; ForceText = $false }
    @{ Row = 14; Col = 5; Value = "  +0.69%  "; ForceText = $false }
    @{ Row = 15; Col = 5; Value = "  +2.57%  "; ForceText = $false }
    @{ Row = 16; Col = 4; Value = "65.69"; ForceText = $true }
    @{ Row = 16; Col = 5; Value = "  +1.99%  "; ForceText = $false }
    @{ Row = 17; Col = 4; Value = "27.881.56"; ForceText = $false }
    @{ Row = 17; Col = 5; Value = "  +1.15%  "; ForceText = $false }
    @{ Row = 18; Col = 4; Value = "231.55"; ForceText = $true }
    @{ Row = 18; Col = 5; Value = "  +0.49%  "; ForceText = $false }
    @{ Row = 19; Col = 5; Value = "  +0.34%  "; ForceText = $false }
    @{ Row = 20; Col = 5; Value = "  +1.08%  "; ForceText = $false }
    @{ Row = 22; Col = 4; Value = "10.78"; ForceText = $true }
    @{ Row = 22; Col = 5; Value = "  +8.57%  "; ForceText = $false }
    @{ Row = 24; Col = 4; Value = "2.15"; ForceText = $true }
    @{ Row = 24; Col = 5; Value = "  +4.05%  "; ForceText = $false }
    @{ Row = 25; Col = 4; Value = "151.00"; ForceText = $true }
    @{ Row = 25; Col = 5; Value = "  +1.19%  "; ForceText = $false }
    @{ Row = 26; Col = 4; Value = "6.90"; ForceText = $true }
    @{ Row = 26; Col = 5; Value = "  +0.37%  "; ForceText = $false }
    @{ Row = 27; Col = 4; Value = "0.112"; ForceText = $true }
    @{ Row = 27; Col = 5; Value = "  +0.47%  "; ForceText = $false }
    @{ Row = 28; Col = 5; Value = "  +0.71%  "; ForceText = $false }
    @{ Row = 29; Col = 5; Value = "  -0.16%  "; ForceText = $false }
    @{ Row = 30; Col = 4; Value = "1.19"; ForceText = $true }
    @{ Row = 30; Col = 5; Value = "  +0.62%  "; ForceText = $false }
    @{ Row = 31; Col = 5; Value = "  +0.04%  "; ForceText = $false }
    @{ Row = 32; Col = 5; Value = "  +0.75%  "; ForceText = $false }
    @{ Row = 33; Col = 4; Value = "1.457.08"; ForceText = $false }
    @{ Row = 33; Col = 5; Value = "  +0.06%  "; ForceText = $false }
    @{ Row = 34; Col = 5; Value = "  +0.48%  "; ForceText = $false }
    @{ Row = 35; Col = 5; Value = "  +1.04%  "; ForceText = $false }
    @{ Row = 36; Col = 5; Value = "  -0.53%  "; ForceText = $false }
    @{ Row = 37; Col = 4; Value = "0.889"; ForceText = $true }
    @{ Row = 37; Col = 5; Value = "  +2.69%  "; ForceText = $false }
    @{ Row = 38; Col = 4; Value = "0.564"; ForceText = $true }
    @{ Row = 38; Col = 5; Value = "  +0.13%  "; ForceText = $false }
    @{ Row = 39; Col = 5; Value = "  +0.32%  "; ForceText = $false }
    @{ Row = 40; Col = 4; Value = "0.916"; ForceText = $true }
    @{ Row = 40; Col = 5; Value = "  -3.56%  "; ForceText = $false }
    @{ Row = 41; Col = 4; Value = "69.21"; ForceText = $true }
    @{ Row = 41; Col = 5; Value = "  -0.11%  "; ForceText = $false }
    @{ Row = 42; Col = 5; Value = "  -0.15%  "; ForceText = $false }
    @{ Row = 43; Col = 5; Value = "  +0.43%  "; ForceText = $false }
    @{ Row = 44; Col = 5; Value = "  -0.35%  "; ForceText = $false }
    @{ Row = 45; Col = 5; Value = "  +0.17%  "; ForceText = $false }
    @{ Row = 46; Col = 4; Value = "1.79"; ForceText = $true }
    @{ Row = 46; Col = 5; Value = "  +6.45%  "; ForceText = $false }
    @{ Row = 47; Col = 5; Value = "  -1.69%  "; ForceText = $false }
    @{ Row = 48; Col = 4; Value = "1.780.86"; ForceText = $false }
    @{ Row = 48; Col = 5; Value = "  +0.82%  "; ForceText = $false }
    @{ Row = 49; Col = 4; Value = "88.32"; ForceText = $true }
    @{ Row = 49; Col = 5; Value = "  +2.43%  "; ForceText = $false }
    @{ Row = 50; Col = 4; Value = "0.101"; ForceText = $true }
    @{ Row = 50; Col = 5; Value = "  +2.05%  "; ForceText = $false }
    @{ Row = 51; Col = 2; Value = "BabyDogeCoin"; ForceText = $false }
    @{ Row = 51; Col = 3; Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; ForceText = $false }
    @{ Row = 51; Col = 4; Value = "0.0₆0101"; ForceText = $false }
    @{ Row = 51; Col = 5; Value = "  -3.90%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
